$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 (the "Neutrophils" row) first
$ws.Rows.Item(4).Delete()

# Rename "Inflammatory-Mac" -> "Resolving-Mac" wherever it appears (column D)
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("D3").Value = "Resolving-Mac"

# Update numeric values in row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2046615
$ws.Range("H2").Value = 0.409323
$ws.Range("I2").Value = 0.4984844148647908
$ws.Range("J2").Value = 0.4984844148647908
$ws.Range("M2").Value = 0.04251333333333333
$ws.Range("N2").Value = 0.12754
$ws.Range("Q2").Value = 0.008700842569999999
$ws.Range("R2").Value = 0.05220505541999999
$ws.Range("S2").Value = 0.4984844148647908
$ws.Range("T2").Value = 0.4984844148647908

# Update numeric values in row 3
$ws.Range("G3").Value = 0.205906
$ws.Range("H3").Value = 0.411812
$ws.Range("I3").Value = 0.5015155851352092
$ws.Range("J3").Value = 0.5015155851352092
$ws.Range("M3").Value = 0.04251333333333333
$ws.Range("N3").Value = 0.12754
$ws.Range("Q3").Value = 0.00875375041333333
$ws.Range("R3").Value = 0.05252250247999999
$ws.Range("S3").Value = 0.5015155851352092
$ws.Range("T3").Value = 0.5015155851352092
